$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New list of tickers (replaces the old 30-ticker list with a new 38-ticker list)
$tickers = @(
    "VCIT","ITE","JNK","BNDX","PSK","DJCI","VAW","IAU","IEO","WOOD",
    "VIG","VIGI","DVY","DIV","VYM","SDIV","VWO","VSS","VXUS","VTV",
    "IMTM","MTUM","VBR","VNQ","OHI","VNQI","VDC","VDE","VFH","VHT",
    "IGF","VIS","VGT","XTN","VPU","QAI","SH","OLN"
)

# Clear out any previous content/formatting from the old range before writing the new one
$ws.Cells.Clear() | Out-Null

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}

# Apply an explicit (non-theme) Calibri 11 font to column A - this creates a
# second font/cellXf entry in styles.xml, matching the new style applied to
# every populated cell in the column.
$colA = $ws.Columns.Item(1)
$colA.Font.Name = "Calibri"
$colA.Font.Size = 11

# Approximate the new column width (best achievable value in this runtime).
$colA.ColumnWidth = 7.8

# Update selection to match the new active cell.
$ws.Range("E11").Select()

# Touch page setup so a <pageSetup> element is written for the sheet.
$ws.PageSetup.Orientation = 1
